$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5.856399999999997
$ws.Range("C11").Value = -12.0834
$ws.Range("B12").Value = 4.880999999999997
$ws.Range("B15").Value = 4.714499999999996
$ws.Range("C23").Value = -12.2168
$ws.Range("B27").Value = 6.117800000000004
$ws.Range("B28").Value = 5.015599999999993
$ws.Range("C28").Value = -13.7552
$ws.Range("B31").Value = 5.242599999999999
$ws.Range("B32").Value = 6.864299999999997
$ws.Range("C32").Value = -11.62370000000001
$ws.Range("C34").Value = -11.92990000000001
$ws.Range("B36").Value = 9.501100000000001
$ws.Range("C36").Value = -11.7499
$ws.Range("C37").Value = -12.75219999999999
$ws.Range("B38").Value = 5.2679
$ws.Range("C42").Value = -12.671
$ws.Range("B46").Value = 6.304100000000002
$ws.Range("C49").Value = -13.70779999999999
$ws.Range("B54").Value = 4.977900000000004
$ws.Range("C54").Value = -14.14679999999999
$ws.Range("B55").Value = 4.892099999999997
$ws.Range("B56").Value = 4.573799999999995
$ws.Range("B67").Value = 5.218499999999994
$ws.Range("B69").Value = 5.276699999999996
$ws.Range("B72").Value = 5.198200000000007
$ws.Range("B73").Value = 8.889699999999992
$ws.Range("C78").Value = -12.3347
$ws.Range("C80").Value = -12.19830000000001
$ws.Range("B83").Value = 6.185599999999999
$ws.Range("B86").Value = 5.057600000000003
$ws.Range("B91").Value = 5.167599999999996
$ws.Range("B93").Value = 4.986399999999996
$ws.Range("C97").Value = -12.1801
$ws.Range("B99").Value = 5.962399999999999
$ws.Range("C99").Value = -12.009
$ws.Range("C100").Value = -12.33899999999999
$ws.Range("C101").Value = -12.8897
$ws.Range("B104").Value = 9.751300000000004
$ws.Range("B105").Value = 8.064500000000002
